$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 32 (subject 30): fill in the missing measurement columns C..H
$ws.Range("B32").Value = "F"
$ws.Range("C32").Value = 10
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 10
$ws.Range("F32").Value = 13
$ws.Range("G32").Value = 10
$ws.Range("H32").Value = 10

# Row 33 (subject 31): change condition from G to F, fill measurement columns, update Valid flag
$ws.Range("B33").Value = "F"
$ws.Range("C33").Value = 6
$ws.Range("D33").Value = 12
$ws.Range("E33").Value = 6
$ws.Range("F33").Value = 5
$ws.Range("G33").Value = 3
$ws.Range("H33").Value = 10
$ws.Range("I33").Value = 2

# Update the selection to match the new active cell
$ws.Range("I33").Select()
